$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2006_count")
$ws.Range("B3:D3").ClearContents()
$ws.Range("H3:K3").ClearContents()
$ws.Range("B13:E13").ClearContents()
$ws.Range("H13:K13").ClearContents()
$ws.Rows("53:54").Delete()

$ws = $wb.Worksheets.Item("2006_weighted")
$ws.Range("B3:D3").ClearContents()
$ws.Range("H3:K3").ClearContents()
$ws.Range("B13:E13").ClearContents()
$ws.Range("H13:K13").ClearContents()
$ws.Rows("53:54").Delete()

$ws = $wb.Worksheets.Item("2007_count")
$ws.Range("B3:D3").ClearContents()
$ws.Range("H3:K3").ClearContents()
$ws.Range("B13:D13").ClearContents()
$ws.Range("H13:K13").ClearContents()
$ws.Rows("53:55").Delete()

$ws = $wb.Worksheets.Item("2007_weighted")
$ws.Range("B3:D3").ClearContents()
$ws.Range("H3:K3").ClearContents()
$ws.Range("B13:D13").ClearContents()
$ws.Range("H13:K13").ClearContents()
$ws.Rows("53:55").Delete()

$ws = $wb.Worksheets.Item("2008_count")
$ws.Range("K3").ClearContents()
$ws.Range("B13:D13").ClearContents()
$ws.Range("G13:K13").ClearContents()
$ws.Rows("53:55").Delete()

$ws = $wb.Worksheets.Item("2008_weighted")
$ws.Range("B13:D13").ClearContents()
$ws.Range("G13:K13").ClearContents()
$ws.Rows("53:55").Delete()

$ws = $wb.Worksheets.Item("2009_count")
$ws.Range("K3").ClearContents()
$ws.Range("B13:D13").ClearContents()
$ws.Range("G13:K13").ClearContents()
$ws.Rows("53:55").Delete()

$ws = $wb.Worksheets.Item("2009_weighted")
$ws.Range("K3").ClearContents()
$ws.Range("B13:D13").ClearContents()
$ws.Range("G13:K13").ClearContents()
$ws.Rows("53:55").Delete()

$ws = $wb.Worksheets.Item("2010_count")
$ws.Range("L3").ClearContents()
$ws.Range("B13:D13").ClearContents()
$ws.Range("H13:L13").ClearContents()
$ws.Rows("53:55").Delete()

$ws = $wb.Worksheets.Item("2010_weighted")
$ws.Range("L3").ClearContents()
$ws.Range("B13:D13").ClearContents()
$ws.Range("H13:L13").ClearContents()
$ws.Rows("53:55").Delete()
